$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the descriptive parameter names in column A (rows 17-24) to append
# the short variable-name aliases used elsewhere in the new model script,
# and rewrite row 21's description entirely.
$ws.Range("A17").Value = "taux de diffusion du pathogène par classe L (eta1)"
$ws.Range("A18").Value = "taux de diffusion du pathogène par classe J (eta2)"
$ws.Range("A19").Value = "taux de diffusion du pathogène par classe A (eta3)"
$ws.Range("A20").Value = "taux mortalité P (mpath)"
$ws.Range("A21").Value = "Taux d'infection par les pathogènes libres (trans2)"
$ws.Range("A22").Value = "proportion d'interaction avec l'environnement infectieux de la classe L (a1)"
$ws.Range("A23").Value = "proportion d'interaction avec l'environnement infectieux de la classe J (a2)"
$ws.Range("A24").Value = "proportion d'interaction avec l'environnement infectieux de la classe A (a3)"

# Bring the formatting of B20/B22/B23/B24 in line with the rest of the
# "value" column (same look as B17-B19) instead of their former, slightly
# divergent style.
$ws.Range("B17").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column A now holds longer text -- widen it to fit (mirrors Excel's
# best-fit column-width recompute after the text grew).
$ws.Columns.Item(1).ColumnWidth = 62.6

# Move the cursor/scroll position to just past the last row, as after
# finishing the edits, and reset the frozen/leftmost view to the top.
$ws.Range("A25").Select()
